# Applies updated market-price-driven profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# columns H:N across the 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 6
$ws_ALC.Range("H6").Value = 61.75
$ws_ALC.Range("I6").Value = 61.75
$ws_ALC.Range("K6").Value = 185.25
$ws_ALC.Range("M6").Value = -73.25

# ALC row 40
$ws_ALC.Range("H40").Value = 5250.0586
$ws_ALC.Range("I40").Value = 4176
$ws_ALC.Range("J40").Value = 6458.375
$ws_ALC.Range("K40").Value = 4176
$ws_ALC.Range("L40").Value = 6458.375
$ws_ALC.Range("M40").Value = -4001
$ws_ALC.Range("N40").Value = -6808.375

# ALC row 64
$ws_ALC.Range("H64").Value = 7081.5386
$ws_ALC.Range("I64").Value = 5225.6665
$ws_ALC.Range("J64").Value = 7638.3
$ws_ALC.Range("K64").Value = 5225.6665
$ws_ALC.Range("L64").Value = 7638.3
$ws_ALC.Range("M64").Value = -4977.6665
$ws_ALC.Range("N64").Value = -8134.3

# ALC row 67
$ws_ALC.Range("H67").Value = 7081.5386
$ws_ALC.Range("I67").Value = 5225.6665
$ws_ALC.Range("J67").Value = 7638.3
$ws_ALC.Range("K67").Value = 5225.6665
$ws_ALC.Range("L67").Value = 7638.3
$ws_ALC.Range("M67").Value = -4367.6665
$ws_ALC.Range("N67").Value = -9354.299999999999

# ALC row 80
$ws_ALC.Range("H80").Value = 1418.3077
$ws_ALC.Range("J80").Value = 1498.8
$ws_ALC.Range("L80").Value = 4496.4
$ws_ALC.Range("N80").Value = -6492.4

# ALC row 83
$ws_ALC.Range("H83").Value = 1418.3077
$ws_ALC.Range("J83").Value = 1498.8
$ws_ALC.Range("L83").Value = 13489.2
$ws_ALC.Range("N83").Value = -23473.2

# ALC row 88
$ws_ALC.Range("H88").Value = 866.36365
$ws_ALC.Range("J88").Value = 498
$ws_ALC.Range("L88").Value = 498
$ws_ALC.Range("N88").Value = -1310

# ALC row 91
$ws_ALC.Range("H91").Value = 866.36365
$ws_ALC.Range("J91").Value = 498
$ws_ALC.Range("L91").Value = 498
$ws_ALC.Range("N91").Value = -3306

# ALC row 95
$ws_ALC.Range("H95").Value = 37599.4
$ws_ALC.Range("J95").Value = 37599.4
$ws_ALC.Range("L95").Value = 37599.4
$ws_ALC.Range("N95").Value = -43091.4

# ALC row 138
$ws_ALC.Range("H138").Value = 2646.6365
$ws_ALC.Range("I138").Value = 1954.5714
$ws_ALC.Range("J138").Value = 3857.75
$ws_ALC.Range("K138").Value = 5863.7142
$ws_ALC.Range("L138").Value = 11573.25
$ws_ALC.Range("M138").Value = -723.7142000000003
$ws_ALC.Range("N138").Value = -21853.25

# ARM row 63
$ws_ARM.Range("H63").Value = 1823.25
$ws_ARM.Range("I63").Value = 1823.25
$ws_ARM.Range("K63").Value = 1823.25
$ws_ARM.Range("M63").Value = -1137.25

# ARM row 66
$ws_ARM.Range("H66").Value = 1823.25
$ws_ARM.Range("I66").Value = 1823.25
$ws_ARM.Range("K66").Value = 9116.25
$ws_ARM.Range("M66").Value = -5684.25

# ARM row 88
$ws_ARM.Range("H88").Value = 1854.5454
$ws_ARM.Range("I88").Value = 1770
$ws_ARM.Range("J88").Value = 1925
$ws_ARM.Range("K88").Value = 1770
$ws_ARM.Range("L88").Value = 1925
$ws_ARM.Range("M88").Value = -1364
$ws_ARM.Range("N88").Value = -2737

# ARM row 91
$ws_ARM.Range("H91").Value = 1854.5454
$ws_ARM.Range("I91").Value = 1770
$ws_ARM.Range("J91").Value = 1925
$ws_ARM.Range("K91").Value = 1770
$ws_ARM.Range("L91").Value = 1925
$ws_ARM.Range("M91").Value = -366
$ws_ARM.Range("N91").Value = -4733

# ARM row 121
$ws_ARM.Range("H121").Value = 0
$ws_ARM.Range("J121").Value = 0
$ws_ARM.Range("L121").Value = 0
$ws_ARM.Range("N121").ClearContents()

# ARM row 128
$ws_ARM.Range("H128").Value = 350000
$ws_ARM.Range("J128").Value = 350000
$ws_ARM.Range("L128").Value = 350000
$ws_ARM.Range("N128").Value = -359960

# BSM row 20
$ws_BSM.Range("H20").Value = 4511.4287
$ws_BSM.Range("I20").Value = 4931.6665
$ws_BSM.Range("K20").Value = 4931.6665
$ws_BSM.Range("M20").Value = -4684.6665

# BSM row 86
$ws_BSM.Range("H86").Value = 5220.1904
$ws_BSM.Range("J86").Value = 8185.5713
$ws_BSM.Range("L86").Value = 8185.5713
$ws_BSM.Range("N86").Value = -10431.5713

# BSM row 89
$ws_BSM.Range("H89").Value = 5220.1904
$ws_BSM.Range("J89").Value = 8185.5713
$ws_BSM.Range("L89").Value = 40927.85649999999
$ws_BSM.Range("N89").Value = -52159.85649999999

# BSM row 99
$ws_BSM.Range("H99").Value = 2165.625
$ws_BSM.Range("I99").Value = 2373.7144
$ws_BSM.Range("K99").Value = 2373.7144
$ws_BSM.Range("M99").Value = -875.7143999999998

# BSM row 105
$ws_BSM.Range("H105").Value = 1846.2222
$ws_BSM.Range("I105").Value = 1825.875
$ws_BSM.Range("J105").Value = 2009
$ws_BSM.Range("K105").Value = 1825.875
$ws_BSM.Range("L105").Value = 2009
$ws_BSM.Range("M105").Value = -78.875
$ws_BSM.Range("N105").Value = -5503

# CRP row 58
$ws_CRP.Range("H58").Value = 3431.348
$ws_CRP.Range("I58").Value = 2722.2632
$ws_CRP.Range("J58").Value = 6799.5
$ws_CRP.Range("K58").Value = 2722.2632
$ws_CRP.Range("L58").Value = 6799.5
$ws_CRP.Range("M58").Value = -2519.2632
$ws_CRP.Range("N58").Value = -7205.5

# CRP row 62
$ws_CRP.Range("H62").Value = 3875
$ws_CRP.Range("J62").Value = 5000
$ws_CRP.Range("L62").Value = 5000
$ws_CRP.Range("N62").Value = -6248

# CRP row 65
$ws_CRP.Range("H65").Value = 3875
$ws_CRP.Range("J65").Value = 5000
$ws_CRP.Range("L65").Value = 25000
$ws_CRP.Range("N65").Value = -31240

# CRP row 132
$ws_CRP.Range("H132").Value = 2490.125
$ws_CRP.Range("I132").Value = 1989.8
$ws_CRP.Range("K132").Value = 5969.4
$ws_CRP.Range("M132").Value = -3439.4

# CRP row 134
$ws_CRP.Range("H134").Value = 1971.9231
$ws_CRP.Range("I134").Value = 1679.1351
$ws_CRP.Range("J134").Value = 7388.5
$ws_CRP.Range("K134").Value = 5037.4053
$ws_CRP.Range("L134").Value = 22165.5
$ws_CRP.Range("M134").Value = -2502.4053
$ws_CRP.Range("N134").Value = -27235.5

# CRP row 136
$ws_CRP.Range("H136").Value = 3431.348
$ws_CRP.Range("I136").Value = 2722.2632
$ws_CRP.Range("J136").Value = 6799.5
$ws_CRP.Range("K136").Value = 8166.7896
$ws_CRP.Range("L136").Value = 20398.5
$ws_CRP.Range("M136").Value = -5616.7896
$ws_CRP.Range("N136").Value = -25498.5

# CUL row 129
$ws_CUL.Range("H129").Value = 1340.909
$ws_CUL.Range("I129").Value = 870
$ws_CUL.Range("J129").Value = 1733.3334
$ws_CUL.Range("K129").Value = 2610
$ws_CUL.Range("L129").Value = 5200.0002
$ws_CUL.Range("M129").Value = 2390
$ws_CUL.Range("N129").Value = -15200.0002

# CUL row 131
$ws_CUL.Range("H131").Value = 993.13336
$ws_CUL.Range("I131").Value = 974.25
$ws_CUL.Range("K131").Value = 2922.75
$ws_CUL.Range("M131").Value = 2117.25

# GSM row 70
$ws_GSM.Range("H70").Value = 9999.5
$ws_GSM.Range("I70").Value = 9999
$ws_GSM.Range("K70").Value = 9999
$ws_GSM.Range("M70").Value = -9729

# GSM row 73
$ws_GSM.Range("H73").Value = 9999.5
$ws_GSM.Range("I73").Value = 9999
$ws_GSM.Range("K73").Value = 9999
$ws_GSM.Range("M73").Value = -9063

# GSM row 80
$ws_GSM.Range("H80").Value = 10003
$ws_GSM.Range("I80").Value = 0
$ws_GSM.Range("K80").Value = 0
$ws_GSM.Range("M80").ClearContents()

# GSM row 83
$ws_GSM.Range("H83").Value = 10003
$ws_GSM.Range("I83").Value = 0
$ws_GSM.Range("K83").Value = 0
$ws_GSM.Range("M83").ClearContents()

# LTW row 2
$ws_LTW.Range("H2").Value = 3500
$ws_LTW.Range("I2").Value = 3000
$ws_LTW.Range("K2").Value = 3000
$ws_LTW.Range("M2").Value = -2888

# LTW row 22
$ws_LTW.Range("H22").Value = 964.3570999999999
$ws_LTW.Range("I22").Value = 679.5
$ws_LTW.Range("J22").Value = 1344.1666
$ws_LTW.Range("K22").Value = 679.5
$ws_LTW.Range("L22").Value = 1344.1666
$ws_LTW.Range("M22").Value = -384.5
$ws_LTW.Range("N22").Value = -1934.1666

# LTW row 27
$ws_LTW.Range("H27").Value = 964.3570999999999
$ws_LTW.Range("I27").Value = 679.5
$ws_LTW.Range("J27").Value = 1344.1666
$ws_LTW.Range("K27").Value = 679.5
$ws_LTW.Range("L27").Value = 1344.1666
$ws_LTW.Range("M27").Value = -572.5
$ws_LTW.Range("N27").Value = -1558.1666

# LTW row 46
$ws_LTW.Range("H46").Value = 6240.7144
$ws_LTW.Range("I46").Value = 2855.75
$ws_LTW.Range("J46").Value = 7594.7
$ws_LTW.Range("K46").Value = 2855.75
$ws_LTW.Range("L46").Value = 7594.7
$ws_LTW.Range("M46").Value = -2667.75
$ws_LTW.Range("N46").Value = -7970.7

# LTW row 136
$ws_LTW.Range("I136").Value = 4000
$ws_LTW.Range("K136").Value = 12000
$ws_LTW.Range("M136").Value = -9450

# WVR row 47
$ws_WVR.Range("H47").Value = 22500
$ws_WVR.Range("I47").Value = 15000
$ws_WVR.Range("J47").Value = 30000
$ws_WVR.Range("K47").Value = 15000
$ws_WVR.Range("L47").Value = 30000
$ws_WVR.Range("M47").Value = -14428
$ws_WVR.Range("N47").Value = -31144

# WVR row 101
$ws_WVR.Range("H101").Value = 14999.667
$ws_WVR.Range("J101").Value = 14999.667
$ws_WVR.Range("L101").Value = 14999.667
$ws_WVR.Range("N101").Value = -21489.667

# WVR row 132
$ws_WVR.Range("H132").Value = 3503.535
$ws_WVR.Range("I132").Value = 3295.7058
$ws_WVR.Range("J132").Value = 4288.6665
$ws_WVR.Range("K132").Value = 9887.117400000001
$ws_WVR.Range("L132").Value = 12865.9995
$ws_WVR.Range("M132").Value = -7357.117400000001
$ws_WVR.Range("N132").Value = -17925.9995
